# Fruta / hortaliza, semanal
# Insert a new weekly record as row 4 (shifting the existing rows 4-11 down to 5-12)
# and populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 4; this shifts existing rows 4-11 -> 5-12
$ws.Rows.Item(4).Insert()

# Fill the new row 4 with the shared ("constant") columns, copied from the
# surrounding records for this market/product combination.
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 44561
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101006
$ws.Cells.Item(4, 10).Value = "Breva"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 18000
$ws.Cells.Item(4, 16).Value = 18000
$ws.Cells.Item(4, 17).Value = "$/bandeja 6 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(4, 19).Value = 3000
$ws.Cells.Item(4, 20).Value = 6
